$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Delete footer rows 1727-1731 (removed in the edit)
$ws.Range("A1727:A1731").EntireRow.Delete()

# 2. Update header row labels (Spanish -> English column codes)
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# 3. Title-case Spanish prepositions (de/del/la/las/el/los/y) in state/municipality names
#    plus one spelling fix (MonteMorelos -> Montemorelos)
$ws.Range("B8").Value = "Pabellón De Arteaga"
$ws.Range("B9").Value = "Rincón De Romos"
$ws.Range("B10").Value = "San Francisco De Los Romo"
$ws.Range("B11").Value = "San José De Gracia"
$ws.Range("B32").Value = "Amatenango De La Frontera"
$ws.Range("B35").Value = "Bejucal De Ocampo"
$ws.Range("B37").Value = "Benemérito De Las Américas"
$ws.Range("B43").Value = "Chiapa De Corzo"
$ws.Range("B47").Value = "Comitán De Domínguez"
$ws.Range("B71").Value = "Mazapa De Madero"
$ws.Range("B84").Value = "San Cristóbal De Las Casas"
$ws.Range("B121").Value = "Guadalupe Y Calvo"
$ws.Range("B124").Value = "Hidalgo Del Parral"
$ws.Range("B139").Value = "San Francisco Del Oro"
$ws.Range("B143").Value = "Valle De Zaragoza"
$ws.Range("B161").Value = "San Juan De Sabinas"
$ws.Range("B175").Value = "Villa De Álvarez"
$ws.Range("A177").Value = "Ciudad De México"
$ws.Range("B181").Value = "Cuajimalpa De Morelos"
$ws.Range("B196").Value = "Coneto De Comonfort"
$ws.Range("B209").Value = "Nombre De Dios"
$ws.Range("B213").Value = "Pánuco De Coronado"
$ws.Range("B219").Value = "San Juan Del Río"
$ws.Range("B220").Value = "San Luis Del Cordero"
$ws.Range("B221").Value = "San Pedro Del Gallo"
$ws.Range("A231").Value = "Estado De México"
$ws.Range("B231").Value = "Acambay De Ruíz Castañeda"
$ws.Range("B234").Value = "Almoloya De Alquisiras"
$ws.Range("B235").Value = "Almoloya De Juárez"
$ws.Range("B236").Value = "Almoloya Del Río"
$ws.Range("B242").Value = "Atizapán De Zaragoza"
$ws.Range("B250").Value = "Chapa De Mota"
$ws.Range("B255").Value = "Coacalco De Berriozábal"
$ws.Range("B262").Value = "Ecatepec De Morelos"
$ws.Range("B269").Value = "Ixtapan De La Sal"
$ws.Range("B270").Value = "Ixtapan Del Oro"
$ws.Range("B287").Value = "Naucalpan De Juárez"
$ws.Range("B299").Value = "San Antonio La Isla"
$ws.Range("B300").Value = "San Felipe Del Progreso"
$ws.Range("B301").Value = "San Martín De Las Pirámides"
$ws.Range("B303").Value = "San Simón De Guerrero"
$ws.Range("B305").Value = "Soyaniquilpan De Juárez"
$ws.Range("B315").Value = "Tenango Del Valle"
$ws.Range("B327").Value = "Tlalnepantla De Baz"
$ws.Range("B333").Value = "Valle De Bravo"
$ws.Range("B334").Value = "Valle De Chalco Solidaridad"
$ws.Range("B335").Value = "Villa De Allende"
$ws.Range("B336").Value = "Villa Del Carbón"
$ws.Range("B350").Value = "Apaseo El Alto"
$ws.Range("B351").Value = "Apaseo El Grande"
$ws.Range("B359").Value = "Dolores Hidalgo Cuna De La Independencia Nacional"
$ws.Range("B363").Value = "Jaral Del Progreso"
$ws.Range("B371").Value = "Purísima Del Rincón"
$ws.Range("B375").Value = "San Diego De La Unión"
$ws.Range("B377").Value = "San Francisco Del Rincón"
$ws.Range("B379").Value = "San Luis De La Paz"
$ws.Range("B381").Value = "Santa Cruz De Juventino Rosas"
$ws.Range("B383").Value = "Silao De La Victoria"
$ws.Range("B388").Value = "Valle De Santiago"
$ws.Range("B393").Value = "Acapulco De Juárez"
$ws.Range("B396").Value = "Ajuchitlán Del Progreso"
$ws.Range("B397").Value = "Alcozauca De Guerrero"
$ws.Range("B401").Value = "Atenango Del Río"
$ws.Range("B402").Value = "Atlamajalcingo Del Monte"
$ws.Range("B404").Value = "Atoyac De Álvarez"
$ws.Range("B405").Value = "Ayutla De Los Libres"
$ws.Range("B408").Value = "Buenavista De Cuéllar"
$ws.Range("B409").Value = "Chilapa De Álvarez"
$ws.Range("B410").Value = "Chilpancingo De Los Bravo"
$ws.Range("B411").Value = "Coahuayutla De José María Izazaga"
$ws.Range("B416").Value = "Coyuca De Benítez"
$ws.Range("B417").Value = "Coyuca De Catalán"
$ws.Range("B421").Value = "Cuetzala Del Progreso"
$ws.Range("B422").Value = "Cutzamala De Pinzón"
$ws.Range("B428").Value = "Huitzuco De Los Figueroa"
$ws.Range("B429").Value = "Iguala De La Independencia"
$ws.Range("B431").Value = "Ixcateopan De Cuauhtémoc"
$ws.Range("B432").Value = "Zihuatanejo De Azueta"
$ws.Range("B434").Value = "La Unión De Isidoro Montes De Oca"
$ws.Range("B437").Value = "Mártir De Cuilapan"
$ws.Range("B450").Value = "Taxco De Alarcón"
$ws.Range("B452").Value = "Técpan De Galeana"
$ws.Range("B454").Value = "Tepecoacuilco De Trujano"
$ws.Range("B456").Value = "Tixtla De Guerrero"
$ws.Range("B460").Value = "Tlapa De Comonfort"
$ws.Range("B472").Value = "Agua Blanca De Iturbide"
$ws.Range("B477").Value = "Atotonilco De Tula"
$ws.Range("B478").Value = "Atotonilco El Grande"
$ws.Range("B484").Value = "Cuautepec De Hinojosa"
$ws.Range("B490").Value = "Huasca De Ocampo"
$ws.Range("B494").Value = "Huejutla De Reyes"
$ws.Range("B497").Value = "Jacala De Ledezma"
$ws.Range("B504").Value = "Mineral Del Monte"
$ws.Range("B505").Value = "Mixquiahuala De Juárez"
$ws.Range("B506").Value = "Molango De Escamilla"
$ws.Range("B508").Value = "Omitlán De Juárez"
$ws.Range("B509").Value = "Pachuca De Soto"
$ws.Range("B511").Value = "Progreso De Obregón"
$ws.Range("B516").Value = "Santiago De Anaya"
$ws.Range("B517").Value = "Santiago Tulantepec De Lugo Guerrero"
$ws.Range("B521").Value = "Tenango De Doria"
$ws.Range("B523").Value = "Tepehuacán De Guerrero"
$ws.Range("B524").Value = "Tepeji Del Río De Ocampo"
$ws.Range("B526").Value = "Tezontepec De Aldama"
$ws.Range("B532").Value = "Tula De Allende"
$ws.Range("B533").Value = "Tulancingo De Bravo"
$ws.Range("B534").Value = "Villa De Tezontepec"
$ws.Range("B536").Value = "Zacualtipán De Ángeles"
$ws.Range("B537").Value = "Zapotlán De Juárez"
$ws.Range("B542").Value = "Acatlán De Juárez"
$ws.Range("B543").Value = "Ahualulco De Mercado"
$ws.Range("B548").Value = "Atemajac De Brizuela"
$ws.Range("B551").Value = "Atotonilco El Alto"
$ws.Range("B553").Value = "Autlán De Navarro"
$ws.Range("B559").Value = "Cañadas De Obregón"
$ws.Range("B566").Value = "Concepción De Buenos Aires"
$ws.Range("B567").Value = "Cuautitlán De García Barragán"
$ws.Range("B574").Value = "Encarnación De Díaz"
$ws.Range("B581").Value = "Huejuquilla El Alto"
$ws.Range("B582").Value = "Ixtlahuacán De Los Membrillos"
$ws.Range("B583").Value = "Ixtlahuacán Del Río"
$ws.Range("B587").Value = "Jilotlán De Los Dolores"
$ws.Range("B593").Value = "La Manzanilla De La Paz"
$ws.Range("B594").Value = "Lagos De Moreno"
$ws.Range("B602").Value = "Ojuelos De Jalisco"
$ws.Range("B607").Value = "San Cristóbal De La Barranca"
$ws.Range("B608").Value = "San Diego De Alejandría"
$ws.Range("B610").Value = "San Juan De Los Lagos"
$ws.Range("B612").Value = "San Martín De Bolaños"
$ws.Range("B614").Value = "San Miguel El Alto"
$ws.Range("B615").Value = "Santa María De Los Ángeles"
$ws.Range("B616").Value = "Santa María Del Oro"
$ws.Range("B619").Value = "Talpa De Allende"
$ws.Range("B620").Value = "Tamazula De Gordiano"
$ws.Range("B623").Value = "Techaluta De Montenegro"
$ws.Range("B626").Value = "Teocuitatlán De Corona"
$ws.Range("B627").Value = "Tepatitlán De Morelos"
$ws.Range("B630").Value = "Tizapán El Alto"
$ws.Range("B631").Value = "Tlajomulco De Zúñiga"
$ws.Range("B643").Value = "Unión De San Antonio"
$ws.Range("B644").Value = "Unión De Tula"
$ws.Range("B645").Value = "Valle De Guadalupe"
$ws.Range("B646").Value = "Valle De Juárez"
$ws.Range("B651").Value = "Yahualica De González Gallo"
$ws.Range("B652").Value = "Zacoalco De Torres"
$ws.Range("B655").Value = "Zapotlán Del Rey"
$ws.Range("B656").Value = "Zapotlán El Grande"
$ws.Range("B682").Value = "Coalcomán De Vázquez Pallares"
$ws.Range("B684").Value = "Cojumatlán De Régules"
$ws.Range("B751").Value = "Tiquicheo De Nicolás Romero"
$ws.Range("B777").Value = "Coatlán Del Río"
$ws.Range("B788").Value = "Puente De Ixtla"
$ws.Range("B794").Value = "Tetela Del Volcán"
$ws.Range("B796").Value = "Tlaltizapán De Zapata"
$ws.Range("B803").Value = "Zacualpan De Amilpas"
$ws.Range("B807").Value = "Bahía De Banderas"
$ws.Range("B809").Value = "Ixtlán Del Río"
$ws.Range("B816").Value = "Santa María Del Oro"
$ws.Range("B835").Value = "Mier Y Noriega"
$ws.Range("B836").Value = "Montemorelos"
$ws.Range("B841").Value = "San Nicolás De Los Garza"
$ws.Range("B848").Value = "Acatlán De Pérez Figueroa"
$ws.Range("B854").Value = "Ayoquezco De Aldama"
$ws.Range("B857").Value = "Chalcatongo De Hidalgo"
$ws.Range("B858").Value = "Ciénega De Zimatlán"
$ws.Range("B859").Value = "Coicoyán De Las Flores"
$ws.Range("B863").Value = "Cuilápam De Guerrero"
$ws.Range("B864").Value = "Cuyamecalco Villa De Zaragoza"
$ws.Range("B865").Value = "El Barrio De La Soledad"
$ws.Range("B866").Value = "Guadalupe De Ramírez"
$ws.Range("B868").Value = "Heroica Ciudad De Ejutla De Crespo"
$ws.Range("B869").Value = "Heroica Ciudad De Huajuapan De León"
$ws.Range("B870").Value = "Heroica Ciudad De Tlaxiaco"
$ws.Range("B872").Value = "Huautla De Jiménez"
$ws.Range("B874").Value = "Ixtlán De Juárez"
$ws.Range("B875").Value = "Heroica Ciudad De Juchitán De Zaragoza"
$ws.Range("B886").Value = "Mariscala De Juárez"
$ws.Range("B887").Value = "Mártires De Tacubaya"
$ws.Range("B889").Value = "Mazatlán Villa De Flores"
$ws.Range("B891").Value = "Miahuatlán De Porfirio Díaz"
$ws.Range("B892").Value = "Nejapa De Madero"
$ws.Range("B894").Value = "Oaxaca De Juárez"
$ws.Range("B895").Value = "Ocotlán De Morelos"
$ws.Range("B896").Value = "Pinotepa De Don Luis"
$ws.Range("B898").Value = "Putla Villa De Guerrero"
$ws.Range("B899").Value = "Reforma De Pineda"
$ws.Range("B900").Value = "Rojas De Cuauhtémoc"
$ws.Range("B904").Value = "San Agustín De Las Juntas"
$ws.Range("B922").Value = "San Felipe Jalapa De Díaz"
$ws.Range("B927").Value = "San Francisco Del Mar"
$ws.Range("B947").Value = "San Juan Bautista Lo De Soto"
$ws.Range("B957").Value = "San Juan De Los Cués"
$ws.Range("B990").Value = "San Miguel Del Puerto"
$ws.Range("B991").Value = "San Miguel El Grande"
$ws.Range("B1010").Value = "San Pedro El Alto"
$ws.Range("B1026").Value = "San Pedro Y San Pablo Ayutla"
$ws.Range("B1027").Value = "San Pedro Y San Pablo Tequixtepec"
$ws.Range("B1048").Value = "Santa Inés Del Monte"
$ws.Range("B1049").Value = "Santa Lucía Del Camino"
$ws.Range("B1062").Value = "Santa María Jalapa Del Marqués"
$ws.Range("B1096").Value = "Santo Domingo De Morelos"
$ws.Range("B1108").Value = "Tamazulápam Del Espíritu Santo"
$ws.Range("B1109").Value = "Tataltepec De Valdés"
$ws.Range("B1110").Value = "Teotitlán De Flores Magón"
$ws.Range("B1111").Value = "Tezoatlán De Segura Y Luna"
$ws.Range("B1112").Value = "Tlacolula De Matamoros"
$ws.Range("B1113").Value = "Villa De Etla"
$ws.Range("B1114").Value = "Villa De Tututepec De Melchor Ocampo"
$ws.Range("B1115").Value = "Villa Sola De Vega"
$ws.Range("B1116").Value = "Villa Tejúpam De La Unión"
$ws.Range("B1118").Value = "Zapotitlán Del Río"
$ws.Range("B1120").Value = "Zimatlán De Álvarez"
$ws.Range("B1143").Value = "Ayotoxco De Guerrero"
$ws.Range("B1147").Value = "Chalchicomula De Sesma"
$ws.Range("B1156").Value = "Chila De La Sal"
$ws.Range("B1167").Value = "Cuayuca De Andrade"
$ws.Range("B1168").Value = "Cuetzalan Del Progreso"
$ws.Range("B1183").Value = "Huehuetlán El Chico"
$ws.Range("B1184").Value = "Huehuetlán El Grande"
$ws.Range("B1188").Value = "Ixcamilpa De Guerrero"
$ws.Range("B1190").Value = "Izúcar De Matamoros"
$ws.Range("B1198").Value = "Los Reyes De Juárez"
$ws.Range("B1207").Value = "Palmar De Bravo"
$ws.Range("B1228").Value = "San Nicolás De Los Ranchos"
$ws.Range("B1232").Value = "San Salvador El Seco"
$ws.Range("B1233").Value = "San Salvador El Verde"
$ws.Range("B1238").Value = "Tecali De Herrera"
$ws.Range("B1246").Value = "Tepanco De López"
$ws.Range("B1247").Value = "Tepatlaxco De Hidalgo"
$ws.Range("B1253").Value = "Tepexi De Rodríguez"
$ws.Range("B1254").Value = "Tetela De Ocampo"
$ws.Range("B1259").Value = "Tlacotepec De Benito Juárez"
$ws.Range("B1275").Value = "Xayacatlán De Bravo"
$ws.Range("B1280").Value = "Xochitlán De Vicente Suárez"
$ws.Range("B1292").Value = "Amealco De Bonfil"
$ws.Range("B1294").Value = "Cadereyta De Montes"
$ws.Range("B1299").Value = "Jalpan De Serra"
$ws.Range("B1300").Value = "Landa De Matamoros"
$ws.Range("B1302").Value = "Pinal De Amoles"
$ws.Range("B1305").Value = "San Juan Del Río"
$ws.Range("B1317").Value = "Armadillo De Los Infante"
$ws.Range("B1323").Value = "Ciudad Del Maíz"
$ws.Range("B1331").Value = "Mexquitic De Carmona"
$ws.Range("B1337").Value = "San Ciro De Acosta"
$ws.Range("B1343").Value = "Santa María Del Río"
$ws.Range("B1345").Value = "Soledad De Graciano Sánchez"
$ws.Range("B1351").Value = "Tanquián De Escobedo"
$ws.Range("B1355").Value = "Villa De Arista"
$ws.Range("B1356").Value = "Villa De Arriaga"
$ws.Range("B1357").Value = "Villa De Guadalupe"
$ws.Range("B1358").Value = "Villa De La Paz"
$ws.Range("B1359").Value = "Villa De Ramos"
$ws.Range("B1360").Value = "Villa De Reyes"
$ws.Range("B1399").Value = "Nacozari De García"
$ws.Range("B1439").Value = "Soto La Marina"
$ws.Range("B1446").Value = "Acuamanala De Miguel Hidalgo"
$ws.Range("B1448").Value = "Amaxac De Guerrero"
$ws.Range("B1454").Value = "Contla De Juan Cuamatzi"
$ws.Range("B1461").Value = "Muñoz De Domingo Arenas"
$ws.Range("B1462").Value = "Nanacamilpa De Mariano Arista"
$ws.Range("B1465").Value = "Papalotla De Xicohténcatl"
$ws.Range("B1466").Value = "San Pablo Del Monte"
$ws.Range("B1471").Value = "Tepetitla De Lardizábal"
$ws.Range("B1474").Value = "Tetla De La Solidaridad"
$ws.Range("B1486").Value = "Ziltlaltépec De Trinidad Sánchez Santos"
$ws.Range("B1494").Value = "Alto Lucero De Gutiérrez Barrios"
$ws.Range("B1497").Value = "Amatlán De Los Reyes"
$ws.Range("B1506").Value = "Boca Del Río"
$ws.Range("B1507").Value = "Camarón De Tejeda"
$ws.Range("B1510").Value = "Castillo De Teayo"
$ws.Range("B1519").Value = "Chinampa De Gorostiza"
$ws.Range("B1529").Value = "Cosamaloapan De Carpio"
$ws.Range("B1530").Value = "Cosautlán De Carvajal"
$ws.Range("B1544").Value = "Hueyapan De Ocampo"
$ws.Range("B1545").Value = "Huiloapan De Cuauhtémoc"
$ws.Range("B1546").Value = "Ignacio De La Llave"
$ws.Range("B1549").Value = "Ixhuacán De Los Reyes"
$ws.Range("B1550").Value = "Ixhuatlán De Madero"
$ws.Range("B1551").Value = "Ixhuatlán Del Café"
$ws.Range("B1552").Value = "Ixhuatlán Del Sureste"
$ws.Range("B1562").Value = "Juchique De Ferrer"
$ws.Range("B1565").Value = "Landero Y Coss"
$ws.Range("B1568").Value = "Las Vigas De Ramírez"
$ws.Range("B1569").Value = "Lerdo De Tejada"
$ws.Range("B1573").Value = "Martínez De La Torre"
$ws.Range("B1574").Value = "Medellín De Bravo"
$ws.Range("B1578").Value = "Mixtla De Altamirano"
$ws.Range("B1580").Value = "Nanchital De Lázaro Cárdenas Del Río"
$ws.Range("B1590").Value = "Ozuluama De Mascareñas"
$ws.Range("B1593").Value = "Paso De Ovejas"
$ws.Range("B1594").Value = "Paso Del Macho"
$ws.Range("B1597").Value = "Poza Rica De Hidalgo"
$ws.Range("B1606").Value = "Sayula De Alemán"
$ws.Range("B1610").Value = "Soledad De Doblado"
$ws.Range("B1616").Value = "Tatahuicapan De Juárez"
$ws.Range("B1649").Value = "Vega De Alatorre"
$ws.Range("B1659").Value = "Zontecomatlán De López Y Fuentes"
$ws.Range("B1660").Value = "Zozocolco De Hidalgo"
$ws.Range("B1677").Value = "Cañitas De Felipe Pescador"
$ws.Range("B1680").Value = "El Plateado De Joaquín Amaro"
$ws.Range("B1690").Value = "Jiménez Del Teul"
$ws.Range("B1696").Value = "Mezquital Del Oro"
$ws.Range("B1701").Value = "Moyahua De Estrada"
$ws.Range("B1702").Value = "Nochistlán De Mejía"
$ws.Range("B1703").Value = "Noria De Ángeles"
$ws.Range("B1713").Value = "Teúl De González Ortega"
$ws.Range("B1714").Value = "Tlaltenango De Sánchez Román"
$ws.Range("B1715").Value = "Trinidad García De La Cadena"
$ws.Range("B1718").Value = "Villa De Cos"

# 4. Precision-adjusted percentage values (column D) - literal replacement to match target bit pattern
$ws.Range("D182").Value = 0.009711546408794116
$ws.Range("D249").Value = 0.0009298289114802876
$ws.Range("D315").Value = 0.0009298289114802876
$ws.Range("D392").Value = 0.0943673030828994
$ws.Range("D395").Value = 0.0009091660467807255
$ws.Range("D401").Value = 0.0009091660467807255
$ws.Range("D418").Value = 0.0009091660467807255
$ws.Range("D484").Value = 0.0009298289114802876
$ws.Range("D775").Value = 0.0009091660467807255
$ws.Range("D819").Value = 0.0009298289114802876
$ws.Range("D898").Value = 0.0009298289114802876
$ws.Range("D1502").Value = 0.0009091660467807255
$ws.Range("D1565").Value = 0.0009298289114802876
$ws.Range("D1587").Value = 0.0009711546408794116
$ws.Range("D1706").Value = 0.0009091660467807255
